# Remove the Guadalajara vs Pachuca match (original row 2, Mexican Liga MX 20:00:00).
# All following rows shift up by one: old row3->2, old row4->3, old row5->4.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(2).Delete()

# Row 2: update odds that differ from the shifted-up source row
$ws.Range("F2").Value = 17.5
$ws.Range("G2").Value = 20
$ws.Range("H2").Value = 1.2
$ws.Range("I2").Value = 1.22
$ws.Range("J2").Value = 7.8
$ws.Range("K2").Value = 8.6
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 1.05
$ws.Range("Q2").Value = 1.01
$ws.Range("R2").Value = 4.7
$ws.Range("S2").Value = 1.06
$ws.Range("T2").Value = 1.01
$ws.Range("U2").Value = 1.01
$ws.Range("V2").Value = 5.5
$ws.Range("W2").Value = 1.05
$ws.Range("X2").Value = 1000
$ws.Range("Y2").Value = 1000
$ws.Range("Z2").Value = 1000
$ws.Range("AA2").Value = 1000
$ws.Range("AB2").Value = 1000
$ws.Range("AC2").Value = 1000
$ws.Range("AD2").Value = 12.5
$ws.Range("AE2").Value = 11
$ws.Range("AF2").Value = 1000
$ws.Range("AG2").Value = 1000
$ws.Range("AH2").Value = 9
$ws.Range("AI2").Value = 8.2
$ws.Range("AJ2").Value = 1000
$ws.Range("AK2").Value = 1000
$ws.Range("AL2").Value = 16.5
$ws.Range("AM2").Value = 13.5
$ws.Range("AN2").Value = 14
$ws.Range("AO2").Value = 4.3

# Row 3: update odds that differ from the shifted-up source row
$ws.Range("F3").Value = 6
$ws.Range("G3").Value = 6.6
$ws.Range("H3").Value = 1.84
$ws.Range("J3").Value = 3.15
$ws.Range("K3").Value = 3.3
$ws.Range("L3").Value = 2.48
$ws.Range("M3").Value = 1.19
$ws.Range("N3").Value = 2.2
$ws.Range("O3").Value = 1.8
$ws.Range("P3").Value = 1.4
$ws.Range("Q3").Value = 3.4
$ws.Range("R3").Value = 1.11
$ws.Range("S3").Value = 8.8
$ws.Range("T3").Value = 2.84
$ws.Range("U3").Value = 1.48
$ws.Range("W3").Value = 1.18
$ws.Range("X3").Value = 6.8
$ws.Range("Y3").Value = 5.1
$ws.Range("Z3").Value = 8.8
$ws.Range("AA3").Value = 23
$ws.Range("AB3").Value = 12.5
$ws.Range("AC3").Value = 8.2
$ws.Range("AD3").Value = 13.5
$ws.Range("AE3").Value = 36
$ws.Range("AF3").Value = 48
$ws.Range("AG3").Value = 32
$ws.Range("AH3").Value = 46
$ws.Range("AI3").Value = 120
$ws.Range("AJ3").Value = 260
$ws.Range("AK3").Value = 190
$ws.Range("AL3").Value = 260
$ws.Range("AM3").Value = 1000
$ws.Range("AN3").Value = 510
$ws.Range("AO3").Value = 42

# Row 4: update odds that differ from the shifted-up source row
$ws.Range("F4").Value = 2.74
$ws.Range("G4").Value = 2.76
$ws.Range("H4").Value = 2.58
$ws.Range("I4").Value = 2.62
$ws.Range("J4").Value = 3.9
$ws.Range("K4").Value = 4
$ws.Range("L4").Value = 1.29
$ws.Range("N4").Value = 6
$ws.Range("P4").Value = 2.66
$ws.Range("Q4").Value = 1.58
$ws.Range("R4").Value = 1.65
$ws.Range("T4").Value = 1.53
$ws.Range("U4").Value = 2.8
$ws.Range("V4").Value = 1.62
$ws.Range("W4").Value = 1.56
$ws.Range("Y4").Value = 16.5
$ws.Range("Z4").Value = 21
$ws.Range("AA4").Value = 36
$ws.Range("AC4").Value = 9.2
$ws.Range("AD4").Value = 12.5
$ws.Range("AE4").Value = 24
$ws.Range("AF4").Value = 22
$ws.Range("AG4").Value = 12.5
$ws.Range("AI4").Value = 28
$ws.Range("AJ4").Value = 42
$ws.Range("AK4").Value = 24
$ws.Range("AN4").Value = 15
$ws.Range("AO4").Value = 13.5

